$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.293.15'
$ws.Range('E2').Value = '  -1.34%  '
$ws.Range('D3').Value = '3.069.27'
$ws.Range('E3').Value = '  -2.35%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '588.63'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').Value = '152.82'
$ws.Range('E6').Value = '  +4.86%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').Value = '0.541'
$ws.Range('E8').Value = '  +2.08%  '
$ws.Range('D9').Value = '3.065.46'
$ws.Range('E9').Value = '  -2.21%  '
$ws.Range('E10').Value = '  -2.57%  '
$ws.Range('E11').Value = '  -3.02%  '
$ws.Range('E12').Value = '  -0.86%  '
$ws.Range('D13').Value = '37.14'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('E14').Value = '  -3.31%  '
$ws.Range('D16').Value = '3.577.97'
$ws.Range('E16').Value = '  -2.33%  '
$ws.Range('D17').Value = '7.12'
$ws.Range('E17').Value = '  -2.59%  '
$ws.Range('D18').Value = '63.339.86'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').Value = '3.069.69'
$ws.Range('E19').Value = '  -2.37%  '
$ws.Range('D20').Value = '473.39'
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('E21').Value = '  +1.04%  '
$ws.Range('E22').Value = '  -3.12%  '
$ws.Range('D23').Value = '7.49'
$ws.Range('E23').Value = '  -1.43%  '
$ws.Range('E24').Value = '  +0.96%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '80.76'
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '12.81'
$ws.Range('E26').Value = '  -1.33%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('E29').Value = '  -2.04%  '
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('E31').Value = '  -2.13%  '
$ws.Range('E32').Value = '  -3.26%  '
$ws.Range('E33').Value = '  +3.54%  '
$ws.Range('D34').Value = '27.04'
$ws.Range('E34').Value = '  -2.32%  '
$ws.Range('E35').Value = '  -1.59%  '
$ws.Range('E36').Value = '  -2.18%  '
$ws.Range('D37').Value = '3.35'
$ws.Range('E37').Value = '  +4.65%  '
$ws.Range('E38').Value = '  -2.15%  '
$ws.Range('E39').Value = '  -4.49%  '
$ws.Range('D40').Value = '50.53'
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('E41').Value = '  -1.33%  '
$ws.Range('D42').Value = '441.23'
$ws.Range('E42').Value = '  -2.83%  '
$ws.Range('D43').Value = '0.282'
$ws.Range('E43').Value = '  -3.59%  '
$ws.Range('E44').Value = '  -3.51%  '
$ws.Range('D45').Value = '39.80'
$ws.Range('E45').Value = '  -1.75%  '
$ws.Range('E46').Value = '  +2.21%  '
$ws.Range('D47').Value = '2.784.45'
$ws.Range('D48').Value = '130.74'
$ws.Range('E48').Value = '  -2.69%  '
$ws.Range('D49').Value = '25.57'
$ws.Range('E49').Value = '  +5.65%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('E51').Value = '  +0.02%  '
